# Append new trading-log rows (104-111) to Sheet1, mirroring the source
# CSV/DB export that produced the earlier rows in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Timestamp, $Action, $Token, $SignalType, $Price, $PositionSizeUsd, $Leverage, $Stiffness, $PnlPercent, $ExitReason, $Status, $ErrorMessage)

    $ws.Cells.Item($Row, 1).Value = $Timestamp
    $ws.Cells.Item($Row, 2).Value = $Action
    $ws.Cells.Item($Row, 3).Value = $Token
    $ws.Cells.Item($Row, 4).Value = $SignalType
    if ($null -ne $Price) { $ws.Cells.Item($Row, 5).Value = $Price }
    if ($null -ne $PositionSizeUsd) { $ws.Cells.Item($Row, 6).Value = $PositionSizeUsd }
    if ($null -ne $Leverage) { $ws.Cells.Item($Row, 7).Value = $Leverage }
    if ($null -ne $Stiffness) { $ws.Cells.Item($Row, 8).Value = $Stiffness }
    if ($null -ne $PnlPercent) { $ws.Cells.Item($Row, 9).Value = $PnlPercent }
    if ($null -ne $ExitReason) { $ws.Cells.Item($Row, 10).Value = $ExitReason }
    $ws.Cells.Item($Row, 11).Value = $Status
    if ($ErrorMessage) { $ws.Cells.Item($Row, 12).Value = $ErrorMessage }
}

Set-Row 104 "2025-11-10T01:44:01.184016" "TRADING_ATTEMPT" "BTC" "UNKNOWN" 105810.2124022901 $null $null $null $null $null "ATTEMPT" "Attempting trade 1/4"

Set-Row 105 "2025-11-10T01:44:02.527362" "POSITION_OPENED" "BTC" "UNKNOWN" 105810.2124022901 3600 40 0.5304736892722409 $null $null "SUCCESS" $null

Set-Row 106 "2025-11-10T01:44:02.565645" "TRADING_ATTEMPT" "NEAR" "UNKNOWN" 2.915819709647294 $null $null $null $null $null "ATTEMPT" "Attempting trade 2/4"

Set-Row 107 "2025-11-10T01:44:03.891119" "POSITION_FAILED" "NEAR" "UNKNOWN" $null $null $null $null $null $null "FAILED" "Trade execution failed for trade 2"

Set-Row 108 "2025-11-10T01:44:03.942391" "TRADING_ATTEMPT" "XRP" "UNKNOWN" 2.401335169398255 $null $null $null $null $null "ATTEMPT" "Attempting trade 3/4"

Set-Row 109 "2025-11-10T01:44:05.108886" "POSITION_OPENED" "XRP" "UNKNOWN" 2.401335169398255 1800 20 0.3424941803961721 $null $null "SUCCESS" $null

Set-Row 110 "2025-11-10T01:44:05.148365" "TRADING_ATTEMPT" "ENA" "UNKNOWN" 0.3440884159654859 $null $null $null $null $null "ATTEMPT" "Attempting trade 4/4"

Set-Row 111 "2025-11-10T01:44:06.376962" "POSITION_OPENED" "ENA" "UNKNOWN" 0.3440884159654859 900 10 0.4929621216371939 $null $null "SUCCESS" $null
